$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New/changed Cypher query text blocks (here-strings avoid any escaping issues)
$newCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
WHERE diag.primary_disease_site IN ['Lymph Node']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age
RETURN  
       coalesce(c.case_id, '') AS `Case ID`,
       coalesce(s.clinical_study_designation, '') AS `Study Code`,
       coalesce(s.clinical_study_type, '') AS  `Study Type`,
       coalesce(demo.breed, '') AS Breed ,
       coalesce(diag.disease_term, '') AS Diagnosis ,
       coalesce(diag.stage_of_disease, '') AS `Stage of Disease`,
       CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END AS Age,
       coalesce(demo.sex, '') AS Sex,
       coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
    coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS `Weight (kg)`,
       coalesce(diag.best_response, '') AS `Response to Treatment`,
       coalesce(co.cohort_description, '') AS `Cohort`
'@
$newStatsQueryDup = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE diag.primary_disease_site IN ['Lymph Node']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@
$newSamplesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE diag.primary_disease_site IN ['Lymph Node']
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed,
        coalesce(diag.disease_term,'') AS Diagnosis, 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
'@
$origStatsQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE diag.primary_disease_site IN ['Lymph Node']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@
$newFilesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f)-[*]->(samp:sample)
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE diag.primary_disease_site IN ['Lymph Node']
OPTIONAL MATCH (s:study)<--(c)<--(diag:diagnosis)<-[*]-(samp:sample)
WITH
        f, parent, c, demo, diag, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value, 
        10^precision AS factor,
        units[i] as unit
WITH    
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN 
        coalesce(f.file_name, '') AS `File Name`,
        coalesce(f.file_type, '') AS `File Type`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(samp.sample_id, '') AS `Sample ID`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis
'@
$newStudyFilesQuery = @'
MATCH (f:file)-->(s:study)
MATCH (s)<--(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
WHERE diag.primary_disease_site IN ['Lymph Node']
WITH DISTINCT f,  s, c, demo, diag
WITH
        f, c, demo, diag, s,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH    
        f, c, demo, diag, s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
WITH    
        f,  c, demo, diag, s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '') AS `File Name`,
  coalesce(f.file_type, '') AS `File Type`,
  coalesce("study", '') AS `Association`,
  coalesce(f.file_description, '') AS `Description`,
  coalesce(f.file_format, '') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
  coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

# Row 2 (CasesTab): update B2 (query) and C2 (StatQuery)
$ws.Range("B2").Value = $newCasesQuery
$ws.Range("C2").Value = $newStatsQueryDup
$ws.Rows.Item(2).RowHeight = 300

# Row 3 (SamplesTab): update B3 (query) and C3 (StatQuery)
$ws.Range("B3").Value = $newSamplesQuery
$ws.Range("C3").Value = $origStatsQuery

# Row 4 (FilesTab): update B4 (query) and C4 (StatQuery)
$ws.Range("B4").Value = $newFilesQuery
$ws.Range("C4").Value = $origStatsQuery
$ws.Rows.Item(4).RowHeight = 409.5

# Row 5 (new StudyFilesTab row)
$ws.Range("A5").Value = "StudyFilesTab"
$ws.Range("B5").Value = $newStudyFilesQuery
$ws.Range("C5").Value = $origStatsQuery
$ws.Range("D5").Value = $ws.Range("D4").Value2
$ws.Range("E5").Value = $ws.Range("E4").Value2
$ws.Range("B5:C5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 375

# View: zoom to 80% and select C4 as the active cell
$excel.ActiveWindow.Zoom = 80
$ws.Range("C4").Select()
